# This script applies the updated evaluation statistics to the sheet.
# (317 cell values changed across rows 4-13, columns B-AO; see commit "updated results and code")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 0.292
$ws.Range("E4").Value = 0.178
$ws.Range("G4").Value = 0.169
$ws.Range("H4").Value = 0.213
$ws.Range("I4").Value = 0.025
$ws.Range("J4").Value = 0.159
$ws.Range("K4").Value = 0.337
$ws.Range("L4").Value = 0.098
$ws.Range("M4").Value = 0.313
$ws.Range("N4").Value = 0.265
$ws.Range("P4").Value = 0.144
$ws.Range("Q4").Value = 0.524
$ws.Range("R4").Value = 0.219
$ws.Range("S4").Value = 0.468
$ws.Range("T4").Value = 0.268
$ws.Range("V4").Value = 0.292
$ws.Range("W4").Value = 0.257
$ws.Range("X4").Value = 0.042
$ws.Range("Y4").Value = 0.205
$ws.Range("Z4").Value = 0.442
$ws.Range("AA4").Value = 0.125
$ws.Range("AB4").Value = 0.354
$ws.Range("AC4").Value = 0.124
$ws.Range("AE4").Value = 0.079
$ws.Range("AF4").Value = 0.703
$ws.Range("AI4").Value = 0.6830000000000001
$ws.Range("AJ4").Value = 0.156
$ws.Range("AK4").Value = 0.394
$ws.Range("AL4").Value = 0.704
$ws.Range("AM4").Value = 0.116
$ws.Range("AN4").Value = 0.341
$ws.Range("AO4").Value = 0.697

# Row 5
$ws.Range("B5").Value = 0.829
$ws.Range("C5").Value = 0.142
$ws.Range("D5").Value = 0.376
$ws.Range("E5").Value = 0.707
$ws.Range("F5").Value = 0.207
$ws.Range("G5").Value = 0.455
$ws.Range("H5").Value = 0.854
$ws.Range("I5").Value = 0.125
$ws.Range("J5").Value = 0.353
$ws.Range("K5").Value = 0.659
$ws.Range("L5").Value = 0.225
$ws.Range("M5").Value = 0.474
$ws.Range("N5").Value = 0.829
$ws.Range("O5").Value = 0.142
$ws.Range("P5").Value = 0.376
$ws.Range("Q5").Value = 0.585
$ws.Range("R5").Value = 0.243
$ws.Range("S5").Value = 0.493
$ws.Range("T5").Value = 0.5610000000000001
$ws.Range("U5").Value = 0.246
$ws.Range("V5").Value = 0.496
$ws.Range("W5").Value = 0.78
$ws.Range("X5").Value = 0.171
$ws.Range("Y5").Value = 0.414
$ws.Range("Z5").Value = 0.829
$ws.Range("AA5").Value = 0.142
$ws.Range("AB5").Value = 0.376
$ws.Range("AC5").Value = 0.756
$ws.Range("AD5").Value = 0.184
$ws.Range("AE5").Value = 0.429
$ws.Range("AF5").Value = 0.951
$ws.Range("AG5").Value = 0.046
$ws.Range("AH5").Value = 0.215
$ws.Range("AI5").Value = 0.805
$ws.Range("AJ5").Value = 0.157
$ws.Range("AK5").Value = 0.396
$ws.Range("AL5").Value = 0.927
$ws.Range("AM5").Value = 0.068
$ws.Range("AN5").Value = 0.26
$ws.Range("AO5").Value = 0.894

# Row 6
$ws.Range("B6").Value = 0.432
$ws.Range("E6").Value = 0.284
$ws.Range("H6").Value = 0.341
$ws.Range("K6").Value = 0.446
$ws.Range("N6").Value = 0.402
$ws.Range("Q6").Value = 0.553
$ws.Range("T6").Value = 0.363
$ws.Range("W6").Value = 0.387
$ws.Range("Z6").Value = 0.577
$ws.Range("AC6").Value = 0.213
$ws.Range("AF6").Value = 0.8080000000000001
$ws.Range("AI6").Value = 0.739
$ws.Range("AL6").Value = 0.8
$ws.Range("AO6").Value = 0.782

# Row 7
$ws.Range("B7").Value = 0.606
$ws.Range("E7").Value = 0.443
$ws.Range("H7").Value = 0.533
$ws.Range("K7").Value = 0.553
$ws.Range("N7").Value = 0.581
$ws.Range("Q7").Value = 0.572
$ws.Range("T7").Value = 0.46
$ws.Range("W7").Value = 0.554
$ws.Range("Z7").Value = 0.705
$ws.Range("AC7").Value = 0.374
$ws.Range("AF7").Value = 0.888
$ws.Range("AI7").Value = 0.777
$ws.Range("AL7").Value = 0.872
$ws.Range("AO7").Value = 0.846

# Row 8
$ws.Range("B8").Value = 0.75
$ws.Range("C8").Value = 0.146
$ws.Range("D8").Value = 0.382
$ws.Range("E8").Value = 0.602
$ws.Range("F8").Value = 0.186
$ws.Range("H8").Value = 0.749
$ws.Range("I8").Value = 0.135
$ws.Range("J8").Value = 0.368
$ws.Range("K8").Value = 0.583
$ws.Range("M8").Value = 0.449
$ws.Range("N8").Value = 0.746
$ws.Range("O8").Value = 0.143
$ws.Range("P8").Value = 0.378
$ws.Range("Q8").Value = 0.5580000000000001
$ws.Range("R8").Value = 0.23
$ws.Range("S8").Value = 0.479
$ws.Range("T8").Value = 0.488
$ws.Range("U8").Value = 0.209
$ws.Range("V8").Value = 0.457
$ws.Range("W8").Value = 0.702
$ws.Range("X8").Value = 0.163
$ws.Range("Y8").Value = 0.404
$ws.Range("Z8").Value = 0.769
$ws.Range("AA8").Value = 0.143
$ws.Range("AB8").Value = 0.378
$ws.Range("AC8").Value = 0.63
$ws.Range("AE8").Value = 0.42
$ws.Range("AF8").Value = 0.876
$ws.Range("AG8").Value = 0.063
$ws.Range("AH8").Value = 0.251
$ws.Range("AI8").Value = 0.796
$ws.Range("AJ8").Value = 0.157
$ws.Range("AK8").Value = 0.396
$ws.Range("AL8").Value = 0.891
$ws.Range("AM8").Value = 0.075
$ws.Range("AN8").Value = 0.273
$ws.Range("AO8").Value = 0.854

# Row 9
$ws.Range("B9").Value = 0.659
$ws.Range("C9").Value = 0.225
$ws.Range("D9").Value = 0.474
$ws.Range("E9").Value = 0.488
$ws.Range("F9").Value = 0.25
$ws.Range("G9").Value = 0.5
$ws.Range("H9").Value = 0.634
$ws.Range("I9").Value = 0.232
$ws.Range("J9").Value = 0.482
$ws.Range("K9").Value = 0.488
$ws.Range("L9").Value = 0.25
$ws.Range("M9").Value = 0.5
$ws.Range("N9").Value = 0.634
$ws.Range("O9").Value = 0.232
$ws.Range("P9").Value = 0.482
$ws.Range("Q9").Value = 0.512
$ws.Range("R9").Value = 0.25
$ws.Range("S9").Value = 0.5
$ws.Range("T9").Value = 0.39
$ws.Range("U9").Value = 0.238
$ws.Range("V9").Value = 0.488
$ws.Range("W9").Value = 0.585
$ws.Range("X9").Value = 0.243
$ws.Range("Y9").Value = 0.493
$ws.Range("Z9").Value = 0.6830000000000001
$ws.Range("AA9").Value = 0.217
$ws.Range("AB9").Value = 0.465
$ws.Range("AC9").Value = 0.512
$ws.Range("AF9").Value = 0.756
$ws.Range("AG9").Value = 0.184
$ws.Range("AH9").Value = 0.429
$ws.Range("AI9").Value = 0.78
$ws.Range("AJ9").Value = 0.171
$ws.Range("AK9").Value = 0.414
$ws.Range("AL9").Value = 0.829
$ws.Range("AM9").Value = 0.142
$ws.Range("AN9").Value = 0.376
$ws.Range("AO9").Value = 0.788

# Row 10
$ws.Range("B10").Value = 0.78
$ws.Range("C10").Value = 0.171
$ws.Range("D10").Value = 0.414
$ws.Range("E10").Value = 0.634
$ws.Range("F10").Value = 0.232
$ws.Range("G10").Value = 0.482
$ws.Range("H10").Value = 0.78
$ws.Range("I10").Value = 0.171
$ws.Range("J10").Value = 0.414
$ws.Range("K10").Value = 0.659
$ws.Range("L10").Value = 0.225
$ws.Range("M10").Value = 0.474
$ws.Range("N10").Value = 0.805
$ws.Range("O10").Value = 0.157
$ws.Range("P10").Value = 0.396
$ws.Range("Q10").Value = 0.585
$ws.Range("R10").Value = 0.243
$ws.Range("S10").Value = 0.493
$ws.Range("T10").Value = 0.5610000000000001
$ws.Range("U10").Value = 0.246
$ws.Range("V10").Value = 0.496
$ws.Range("W10").Value = 0.78
$ws.Range("X10").Value = 0.171
$ws.Range("Y10").Value = 0.414
$ws.Range("Z10").Value = 0.829
$ws.Range("AA10").Value = 0.142
$ws.Range("AB10").Value = 0.376
$ws.Range("AC10").Value = 0.634
$ws.Range("AD10").Value = 0.232
$ws.Range("AE10").Value = 0.482
$ws.Range("AF10").Value = 0.951
$ws.Range("AG10").Value = 0.046
$ws.Range("AH10").Value = 0.215
$ws.Range("AI10").Value = 0.805
$ws.Range("AJ10").Value = 0.157
$ws.Range("AK10").Value = 0.396
$ws.Range("AL10").Value = 0.927
$ws.Range("AM10").Value = 0.068
$ws.Range("AN10").Value = 0.26
$ws.Range("AO10").Value = 0.894

# Row 11
$ws.Range("B11").Value = 0.829
$ws.Range("C11").Value = 0.142
$ws.Range("D11").Value = 0.376
$ws.Range("E11").Value = 0.707
$ws.Range("F11").Value = 0.207
$ws.Range("G11").Value = 0.455
$ws.Range("H11").Value = 0.854
$ws.Range("I11").Value = 0.125
$ws.Range("J11").Value = 0.353
$ws.Range("K11").Value = 0.659
$ws.Range("L11").Value = 0.225
$ws.Range("M11").Value = 0.474
$ws.Range("N11").Value = 0.829
$ws.Range("O11").Value = 0.142
$ws.Range("P11").Value = 0.376
$ws.Range("Q11").Value = 0.585
$ws.Range("R11").Value = 0.243
$ws.Range("S11").Value = 0.493
$ws.Range("T11").Value = 0.5610000000000001
$ws.Range("U11").Value = 0.246
$ws.Range("V11").Value = 0.496
$ws.Range("W11").Value = 0.78
$ws.Range("X11").Value = 0.171
$ws.Range("Y11").Value = 0.414
$ws.Range("Z11").Value = 0.829
$ws.Range("AA11").Value = 0.142
$ws.Range("AB11").Value = 0.376
$ws.Range("AC11").Value = 0.6830000000000001
$ws.Range("AD11").Value = 0.217
$ws.Range("AE11").Value = 0.465
$ws.Range("AF11").Value = 0.951
$ws.Range("AG11").Value = 0.046
$ws.Range("AH11").Value = 0.215
$ws.Range("AI11").Value = 0.805
$ws.Range("AJ11").Value = 0.157
$ws.Range("AK11").Value = 0.396
$ws.Range("AL11").Value = 0.927
$ws.Range("AM11").Value = 0.068
$ws.Range("AN11").Value = 0.26
$ws.Range("AO11").Value = 0.894

# Row 12
$ws.Range("B12").Value = 1.382
$ws.Range("C12").Value = 0.707
$ws.Range("D12").Value = 0.841
$ws.Range("E12").Value = 1.621
$ws.Range("F12").Value = 1.063
$ws.Range("G12").Value = 1.031
$ws.Range("H12").Value = 1.543
$ws.Range("I12").Value = 1.22
$ws.Range("J12").Value = 1.104
$ws.Range("K12").Value = 1.407
$ws.Range("L12").Value = 0.538
$ws.Range("M12").Value = 0.733
$ws.Range("N12").Value = 1.353
$ws.Range("O12").Value = 0.522
$ws.Range("P12").Value = 0.723
$ws.Range("Z12").Value = 1.235
$ws.Range("AA12").Value = 0.298
$ws.Range("AB12").Value = 0.546
$ws.Range("AC12").Value = 2.032
$ws.Range("AD12").Value = 3.902
$ws.Range("AE12").Value = 1.975
$ws.Range("AF12").Value = 1.231
$ws.Range("AG12").Value = 0.229
$ws.Range("AH12").Value = 0.478
$ws.Range("AI12").Value = 1.03
$ws.Range("AJ12").Value = 0.029
$ws.Range("AK12").Value = 0.171
$ws.Range("AL12").Value = 1.105
$ws.Range("AM12").Value = 0.094
$ws.Range("AN12").Value = 0.307
$ws.Range("AO12").Value = 1.122

# Row 13
$ws.Range("B13").Value = 3.512
$ws.Range("C13").Value = 1.372
$ws.Range("D13").Value = 1.171
$ws.Range("E13").Value = 4.543
$ws.Range("F13").Value = 0.762
$ws.Range("G13").Value = 0.873
$ws.Range("H13").Value = 4.5
$ws.Range("I13").Value = 0.987
$ws.Range("J13").Value = 0.993
$ws.Range("K13").Value = 2.3
$ws.Range("L13").Value = 0.61
$ws.Range("M13").Value = 0.781
$ws.Range("N13").Value = 3.317
$ws.Range("O13").Value = 0.802
$ws.Range("P13").Value = 0.895
$ws.Range("Z13").Value = 2.85
$ws.Range("AA13").Value = 3.978
$ws.Range("AB13").Value = 1.994
$ws.Range("AC13").Value = 6.244
$ws.Range("AD13").Value = 2.965
$ws.Range("AE13").Value = 1.722
$ws.Range("AF13").Value = 1.659
$ws.Range("AG13").Value = 0.713
$ws.Range("AH13").Value = 0.844
$ws.Range("AI13").Value = 1.244
$ws.Range("AJ13").Value = 0.184
$ws.Range("AK13").Value = 0.429
$ws.Range("AL13").Value = 1.634
$ws.Range("AM13").Value = 0.8169999999999999
$ws.Range("AN13").Value = 0.904
$ws.Range("AO13").Value = 1.512
